# Updates cryptos list (prices + 1h volume deltas) per upstream scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.788.43'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').Value = '1.644.88'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.53%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.00'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.501'
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('E7').Value = '  +0.53%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0629'
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0843'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').Value = '1.868.51'
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('D13').Value = '1.655.42'
$ws.Range('E13').Value = '  +1.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.16'
$ws.Range('E14').Value = '  -1.17%  '
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('E16').Value = '  -2.02%  '
$ws.Range('D17').Value = '26.795.45'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').Value = '0.0₃0737'
$ws.Range('E18').Value = '  -1.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '213.64'
$ws.Range('E19').Value = '  -2.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.37'
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('E22').Value = '  +15.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.26'
$ws.Range('E23').Value = '  -0.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.34'
$ws.Range('E24').Value = '  -2.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.46'
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('E26').Value = '  +0.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.118'
$ws.Range('E27').Value = '  -1.48%  '
$ws.Range('E28').Value = '  -0.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.63'
$ws.Range('E29').Value = '  -1.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0508'
$ws.Range('E30').Value = '  -2.02%  '
$ws.Range('E31').Value = '  +0.64%  '
$ws.Range('E32').Value = '  -2.17%  '
$ws.Range('E33').Value = '  -1.89%  '
$ws.Range('D34').Value = '1.295.23'
$ws.Range('E34').Value = '  +1.55%  '
$ws.Range('E35').Value = '  -0.53%  '
$ws.Range('E36').Value = '  +1.46%  '
$ws.Range('E37').Value = '  -4.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.534'
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.824'
$ws.Range('E39').Value = '  -0.60%  '
$ws.Range('E40').Value = '  +0.49%  '
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('E43').Value = '  -2.12%  '
$ws.Range('D44').Value = '1.795.97'
$ws.Range('E44').Value = '  +0.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.65'
$ws.Range('E45').Value = '  +3.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.55'
$ws.Range('E46').Value = '  -1.83%  '
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0525'
$ws.Range('E48').Value = '  +1.49%  '
$ws.Range('E49').Value = '  -1.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0974'
$ws.Range('E50').Value = '  -0.33%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.01'
$ws.Range('E51').Value = '  +0.65%  '
